$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 0.3145330535152152
$ws.Range("J2").Value = 28.0946659139879
$ws.Range("K2").Value = 11.99
$ws.Range("L2").Value = 38.12
$ws.Range("M2").Value = 387
$ws.Range("N2").Value = 42
$ws.Range("F3").Value = 0
$ws.Range("J3").Value = 169.5756191487949
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 24.99
$ws.Range("M3").Value = 0
$ws.Range("N3").Value = 102
$ws.Range("F4").Value = 0.4717977215893304
$ws.Range("J4").Value = -97.8168688737677
$ws.Range("K4").Value = 16.98
$ws.Range("L4").Value = 35.99
$ws.Range("M4").Value = 768
$ws.Range("N4").Value = 223
$ws.Range("F5").Value = 0
$ws.Range("J5").Value = 342.0548743438145
$ws.Range("K5").Value = 16.73
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = 1.446
$ws.Range("N5").Value = 0
$ws.Range("F6").Value = 0.3715613382899628
$ws.Range("J6").Value = -68.1639571430708
$ws.Range("K6").Value = 19.99
$ws.Range("L6").Value = 53.8
$ws.Range("M6").Value = 84
$ws.Range("N6").Value = 19
$ws.Range("F7").Value = 0.3051963893632593
$ws.Range("J7").Value = 241.4145356880589
$ws.Range("K7").Value = 12.51
$ws.Range("L7").Value = 40.99
$ws.Range("M7").Value = 3.586
$ws.Range("N7").Value = 43
$ws.Range("F8").Value = 0.3642767295597484
$ws.Range("J8").Value = -93.5356860437992
$ws.Range("K8").Value = 14.48
$ws.Range("L8").Value = 39.75
$ws.Range("M8").Value = 281
$ws.Range("N8").Value = 18
$ws.Range("F9").Value = 0.3566666666666667
$ws.Range("J9").Value = 91.75987803982652
$ws.Range("K9").Value = 14.98
$ws.Range("L9").Value = 42
$ws.Range("M9").Value = 7.885
$ws.Range("N9").Value = 55
$ws.Range("F10").Value = 0.3827730169193584
$ws.Range("J10").Value = -90.5026128208562
$ws.Range("K10").Value = 17.42
$ws.Range("L10").Value = 45.51
$ws.Range("M10").Value = 247
$ws.Range("N10").Value = 35
$ws.Range("F11").Value = 0.3678215472235406
$ws.Range("J11").Value = -97.05295277126103
$ws.Range("K11").Value = 15.5
$ws.Range("L11").Value = 42.14
$ws.Range("M11").Value = 671
$ws.Range("N11").Value = 19
$ws.Range("F12").Value = 0
$ws.Range("J12").Value = 691.0690657990839
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 40
$ws.Range("M12").Value = 0
$ws.Range("N12").Value = 42
$ws.Range("F13").Value = 0
$ws.Range("J13").Value = 143.2663437951537
$ws.Range("K13").Value = 19
$ws.Range("L13").Value = 0
$ws.Range("M13").Value = 5.311
$ws.Range("N13").Value = 0
$ws.Range("F14").Value = 0.4862797389290477
$ws.Range("J14").Value = -74.93311631577436
$ws.Range("K14").Value = 69.29
$ws.Range("L14").Value = 142.49
$ws.Range("M14").Value = 313
$ws.Range("N14").Value = 60
$ws.Range("F15").Value = 0.4539119417168198
$ws.Range("J15").Value = 22.96479565824042
$ws.Range("K15").Value = 42.99
$ws.Range("L15").Value = 94.71
$ws.Range("M15").Value = 39
$ws.Range("N15").Value = 1

$ws.Range("J16").ClearContents()

$ws.Range("A5").Select()
